# "update data with resort sheetname"
# The workbook currently has the sheets in the order: 2022-Q2, 总计.
# Re-sort the sheet tabs so that 总计 (the summary/total sheet) comes
# first, followed by 2022-Q2 - i.e. move "总计" to be before "2022-Q2".

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")
$ws2022  = $wb.Worksheets.Item("2022-Q2")

# Move "总计" so that it sits immediately before "2022-Q2" -> new order:
# 总计, 2022-Q2
$wsTotal.Move($ws2022)
